$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# AWS pricing corrections (benchmark-aws-ffmpeg pricing fix)
$ws.Range("B25").Value = 0.07
$ws.Range("B26").Value = 0.14
$ws.Range("B27").Value = 0.28
$ws.Range("B28").Value = 0.56
$ws.Range("B32").Value = 0.105
$ws.Range("B33").Value = 0.21
$ws.Range("B34").Value = 0.42
$ws.Range("B35").Value = 0.84

# Update view / selection state to match the saved workbook snapshot
$ws.Range("F42").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$win.TabRatio = 0.211
